$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 (shifts rows 23..65 down to 24..66,
# and Excel carries the used range/dimension down to row 66 automatically).
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly price record.
$ws.Cells.Item(23, 1).Value = 8
$ws.Cells.Item(23, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(23, 3).Value = "Coquimbo"
$ws.Cells.Item(23, 4).Value = 44581
$ws.Cells.Item(23, 5).Value = 4
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100109
$ws.Cells.Item(23, 8).Value = "Uva"
$ws.Cells.Item(23, 9).Value = 100109001
$ws.Cells.Item(23, 10).Value = "Uva"
$ws.Cells.Item(23, 11).Value = "Superior Seedless"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 440
$ws.Cells.Item(23, 14).Value = 12000
$ws.Cells.Item(23, 15).Value = 13000
$ws.Cells.Item(23, 16).Value = 12500
$ws.Cells.Item(23, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(23, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(23, 19).Value = 694
$ws.Cells.Item(23, 20).Value = 18
